# Listening mode bug fix.
# Applies the ToDos.xlsx edits:
#  - remove two completed/duplicate research rows ("Learn PCB Design", "Read Practical Electronics")
#  - replace "Reset Circuit" row with a new todo "Restore Default option in App"
#  - mark several in-progress items as Closed/Rejected (with one new note)
#  - bump the header date

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# --- Remove obsolete rows (delete bottom-most first so row numbers stay valid) ---
$rowReadPE  = $ws.Cells.Find("Read Practical Electronics").Row
$ws.Rows.Item($rowReadPE).EntireRow.Delete() | Out-Null

$rowLearnPCB = $ws.Cells.Find("Learn PCB Design").Row
$ws.Rows.Item($rowLearnPCB).EntireRow.Delete() | Out-Null

# --- Status updates on existing rows ---
$ws.Cells.Find("Improve streaming performance").Offset(0, 1).Value = "Closed"
$ws.Cells.Find("Stabilize startup / WiFi/ cloud connection").Offset(0, 1).Value = "Closed"

$rowSx = $ws.Cells.Find("SX1509 init sometimes fails").Row
$ws.Cells.Item($rowSx, 3).Value = "Closed"
$ws.Cells.Item($rowSx, 4).Value = "Can't reproduce"
$ws.Cells.Item($rowSx, 4).WrapText = $true

$ws.Cells.Find("Hissing, crackling, high pitch noise come from speakers in rest").Offset(0, 1).Value = "Closed"
$ws.Cells.Find("Hola App is draining battery").Offset(0, 1).Value = "Closed"
$ws.Cells.Find("Keep trying to connect to WiFi when in listening mode").Offset(0, 1).Value = "Rejected"

# --- Replace "Reset Circuit" with the new todo item (status stays Open) ---
$rowReset = $ws.Cells.Find("Reset Circuit").Row
$ws.Cells.Item($rowReset, 2).Value = "Restore Default option in App"

# --- Header date bump ---
$ws.Range("F1").Value = 43174

# --- Shrink the print area to match the now-shorter sheet (2 fewer rows) ---
$wb.Names.Item("Sheet1!Print_Area").RefersTo = "=Sheet1!`$A`$1:`$D`$31"

# --- Leave the cursor where editing finished ---
$ws.Activate() | Out-Null
$ws.Range("C84").Select() | Out-Null
